# Cadm3-Cadm3 LR-pairs sheet: refresh with new TPM-derived expression values.
#
# The sheet is a 5x5 sending-cluster x target-cluster matrix (rows 2-26, 5
# row-blocks of 5 rows each). Because the ligand and receptor are both
# "Cadm3", the "Ligand ..." columns (G:J, keyed by sending cluster A) and the
# "Receptor ..." columns (M:P, keyed by target cluster D) mirror the same
# five per-cluster numbers. Re-running the NATMI pipeline on the new TPM
# matrix changed the average/total expression for four of the five clusters
# (FAPs' numbers happen to be identical under the new TPM), which ripples
# into every derived specificity/edge-weight column.
#
# This script re-derives the whole block from the five updated per-cluster
# (average expression, total expression) pairs, reproducing NATMI's formulas:
#   specificity   = cluster value / sum(all cluster values)
#   edge weight   = ligand value * receptor value
#   edge specificity = edge weight / sum(all edge weights)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending/target cluster order as laid out in the sheet (rows 2-26).
$clusters = @("ECs", "FAPs", "Inflammatory-Mac", "MuSCs", "Resolving-Mac")

# Updated per-cluster ligand/receptor average & total expression values
# (Cadm3 is both the ligand and the receptor, so one table serves both).
$avgExpr = @{
    "ECs"               = 5.342589
    "FAPs"              = 4.362043666666667
    "Inflammatory-Mac"  = 0.8673346666666665
    "MuSCs"             = 1.425086
    "Resolving-Mac"     = 0.03909266666666666
}
$totExpr = @{
    "ECs"               = 16.027767
    "FAPs"              = 13.086131
    "Inflammatory-Mac"  = 2.602004
    "MuSCs"             = 2.850172
    "Resolving-Mac"     = 0.117278
}

$sumAvg = 0.0
$sumTot = 0.0
foreach ($cl in $clusters) {
    $sumAvg += $avgExpr[$cl]
    $sumTot += $totExpr[$cl]
}

$specAvg = @{}
$specTot = @{}
foreach ($cl in $clusters) {
    $specAvg[$cl] = $avgExpr[$cl] / $sumAvg
    $specTot[$cl] = $totExpr[$cl] / $sumTot
}

# First pass: write the per-cluster ligand (G:J) and receptor (M:P) columns,
# and the raw (unnormalised) edge weights (Q:R), for every row.
$row = 2
$edgeAvgW = @{}
$edgeTotW = @{}
foreach ($send in $clusters) {
    foreach ($target in $clusters) {
        $ws.Range("G$row").Value = $avgExpr[$send]
        $ws.Range("H$row").Value = $totExpr[$send]
        $ws.Range("I$row").Value = $specAvg[$send]
        $ws.Range("J$row").Value = $specTot[$send]

        $ws.Range("M$row").Value = $avgExpr[$target]
        $ws.Range("N$row").Value = $totExpr[$target]
        $ws.Range("O$row").Value = $specAvg[$target]
        $ws.Range("P$row").Value = $specTot[$target]

        $q = $avgExpr[$send] * $avgExpr[$target]
        $r = $totExpr[$send] * $totExpr[$target]
        $edgeAvgW[$row] = $q
        $edgeTotW[$row] = $r

        $ws.Range("Q$row").Value = $q
        $ws.Range("R$row").Value = $r

        $row++
    }
}

# Second pass: normalise edge weights into S (avg) / T (total) specificity
# across all 25 sending/target combinations.
$sumQ = 0.0
$sumR = 0.0
foreach ($r in 2..26) {
    $sumQ += $edgeAvgW[$r]
    $sumR += $edgeTotW[$r]
}
foreach ($r in 2..26) {
    $ws.Range("S$r").Value = $edgeAvgW[$r] / $sumQ
    $ws.Range("T$r").Value = $edgeTotW[$r] / $sumR
}
